# Update cryptocurrency price/volume figures (cryptos list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.734.42'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '1.599.04'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = "'211.83"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.13%  '
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('E8').Value = '  +0.36%  '
$ws.Range('E9').Value = '  +0.19%  '
$ws.Range('D10').Value = "'19.65"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.16%  '
$ws.Range('E11').Value = '  +0.83%  '
$ws.Range('D12').Value = '1.824.20'
$ws.Range('E12').Value = '  +0.17%  '
$ws.Range('D13').Value = '1.588.37'
$ws.Range('E13').Value = '  -1.20%  '
$ws.Range('E15').Value = '  +0.51%  '
$ws.Range('D16').Value = "'65.02"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.08%  '
$ws.Range('D17').Value = '0.0₃0740'
$ws.Range('E17').Value = '  -1.15%  '
$ws.Range('E18').Value = '  +0.12%  '
$ws.Range('D19').Value = "'208.65"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.10%  '
$ws.Range('D20').Value = "'7.15"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.82%  '
$ws.Range('E21').Value = '  +0.92%  '
$ws.Range('D22').Value = "'2.22"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.94%  '
$ws.Range('E23').Value = '  +1.10%  '
$ws.Range('D24').Value = "'143.96"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.73%  '
$ws.Range('D25').Value = "'1.00"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('D26').Value = "'7.13"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.42%  '
$ws.Range('E27').Value = '  -0.39%  '
$ws.Range('D28').Value = "'15.37"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.33%  '
$ws.Range('D29').Value = "'0.0508"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.55%  '
$ws.Range('E30').Value = '  +0.31%  '
$ws.Range('E31').Value = '  +1.64%  '
$ws.Range('E32').Value = '  +0.71%  '
$ws.Range('D33').Value = '1.274.28'
$ws.Range('E33').Value = '  -0.69%  '
$ws.Range('E34').Value = '  +1.51%  '
$ws.Range('D35').Value = "'1.23"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +16.25%  '
$ws.Range('E36').Value = '  +0.67%  '
$ws.Range('E37').Value = '  -3.81%  '
$ws.Range('E38').Value = '  -1.21%  '
$ws.Range('D39').Value = "'0.825"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.29%  '
$ws.Range('D40').Value = "'5.48"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.98%  '
$ws.Range('E41').Value = '  +0.21%  '
$ws.Range('E42').Value = '  -0.39%  '
$ws.Range('D43').Value = "'62.63"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.09%  '
$ws.Range('D44').Value = '1.736.20'
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('D45').Value = "'90.43"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.12%  '
$ws.Range('D46').Value = "'1.57"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.72%  '
$ws.Range('D47').Value = "'0.102"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.25%  '
$ws.Range('D48').Value = "'0.0513"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.98%  '
$ws.Range('D49').Value = "'7.52"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.28%  '
$ws.Range('D50').Value = "'1.01"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.15%  '
$ws.Range('E51').Value = '  +1.61%  '
